# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Wed Sep 18 19:24:48 UTC 2024 with GitHub Actions"
# Updates Price (D) / Volume(1h) (E) columns, and for two swapped rows also
# Coin (B) and Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '60.448.22'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.322.91'
$ws.Range("E3").Value = '  -1.20%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.94%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.576'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.22%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.325.78'
$ws.Range("E9").Value = '  -0.88%  '

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.06%  '

# Row 11: Toncoin
$ws.Range("E11").Value = '  +0.11%  '

# Row 12: TRON
$ws.Range("E12").Value = '  -0.59%  '

# Row 13: Cardano
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.49%  '

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.14%  '

# Row 15: WrappedBTC
$ws.Range("D15").Value = '60.527.10'
$ws.Range("E15").Value = '  +0.10%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '2.747.88'
$ws.Range("E16").Value = '  -0.82%  '

# Row 17: ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000133'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.08%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.332.56'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.99%  '

# Row 20: Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.52%  '

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '313.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.41%  '

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.95%  '

# Row 23: Dai
$ws.Range("E23").Value = '  -0.01%  '

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.10%  '

# Row 25: Kaspa
$ws.Range("E25").Value = '  +0.28%  '

# Row 26: Binance-PegBSC-USD
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.18%  '

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.81%  '

# Row 28: Fetch.AI
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.12%  '

# Row 29: SuiNetwork
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.38%  '

# Row 30: Monero
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -1.56%  '

# Row 32: PEPE
$ws.Range("D32").Value = '0.0₃0724'
$ws.Range("E32").Value = '  -1.04%  '

# Row 33: Aptos
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.75%  '

# Row 34: ImmutableX
$ws.Range("B34").Value = 'PolygonEcosystemToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.382'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '

# Row 35: PolygonEcosystemToken
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.36'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.68%  '

# Row 36: EthereumClassic
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.32%  '

# Row 37: USDe
$ws.Range("E37").Value = '  -0.01%  '

# Row 38: FirstDigitalUSD
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '

# Row 39: NEARProtocol
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '

# Row 40: Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '320.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.60%  '

# Row 41: Stacks
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.69%  '

# Row 42: OKB
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '

# Row 43: Aave
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.78%  '

# Row 44: Filecoin
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.92%  '

# Row 45: Stellar
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0939'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.55%  '

# Row 46: InjectiveProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.14%  '

# Row 47: Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.567'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '

# Row 48: Hedera
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0496'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.70%  '

# Row 49: VeChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0215'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.35%  '

# Row 50: BabyDogeCoin
$ws.Range("E50").Value = '  +5.26%  '

# Row 51: WhiteBITCoin
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
